$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "G2" = 27.07648166666667
    "H2" = 81.229445
    "I2" = 0.1943552322922666
    "J2" = 0.1943552322922666
    "M2" = 11.01658666666667
    "N2" = 33.04976
    "O2" = 0.1837409300120545
    "P2" = 0.1837409300120545
    "Q2" = 298.2904069092444
    "R2" = 2684.6136621832
    "S2" = 0.03571101113408996
    "T2" = 0.03571101113408996
    "G3" = 27.07648166666667
    "H3" = 81.229445
    "I3" = 0.1943552322922666
    "J3" = 0.1943552322922666
    "O3" = 0.06509859443116503
    "P3" = 0.06509859443116503
    "Q3" = 105.6829647091594
    "R3" = 951.1466823824348
    "S3" = 0.01265225244256913
    "T3" = 0.01265225244256913
    "G4" = 27.07648166666667
    "H4" = 81.229445
    "I4" = 0.1943552322922666
    "J4" = 0.1943552322922666
    "M4" = 4.001997666666667
    "N4" = 12.005993
    "O4" = 0.06674760480978428
    "P4" = 0.06674760480978428
    "Q4" = 108.3600164515428
    "R4" = 975.240148063885
    "S4" = 0.01297274623775804
    "T4" = 0.01297274623775804
    "G5" = 27.07648166666667
    "H5" = 81.229445
    "I5" = 0.1943552322922666
    "J5" = 0.1943552322922666
    "M5" = 41.03546066666667
    "N5" = 123.106382
    "O5" = 0.6844128707469962
    "P5" = 0.6844128707469963
    "Q5" = 1111.095898424221
    "R5" = 9999.863085817989
    "S5" = 0.1330192224778495
    "T5" = 0.1330192224778495
    "I6" = 0.599012687336886
    "J6" = 0.599012687336886
    "M6" = 11.01658666666667
    "N6" = 33.04976
    "O6" = 0.1837409300120545
    "P6" = 0.1837409300120545
    "Q6" = 919.3461690849956
    "R6" = 8274.115521764959
    "S6" = 0.1100631482602994
    "T6" = 0.1100631482602995
    "I7" = 0.599012687336886
    "J7" = 0.599012687336886
    "O7" = 0.06509859443116503
    "P7" = 0.06509859443116503
    "S7" = 0.0389948839920662
    "T7" = 0.0389948839920662
    "I8" = 0.599012687336886
    "J8" = 0.599012687336886
    "M8" = 4.001997666666667
    "N8" = 12.005993
    "O8" = 0.06674760480978428
    "P8" = 0.06674760480978428
    "Q8" = 333.9710687947893
    "R8" = 3005.739619153103
    "S8" = 0.03998266213040934
    "T8" = 0.03998266213040934
    "I9" = 0.599012687336886
    "J9" = 0.599012687336886
    "M9" = 41.03546066666667
    "N9" = 123.106382
    "O9" = 0.6844128707469962
    "P9" = 0.6844128707469963
    "Q9" = 3424.453934963947
    "R9" = 30820.08541467552
    "S9" = 0.409971992954111
    "T9" = 0.409971992954111
    "G10" = 14.445417
    "H10" = 43.336251
    "I10" = 0.1036893349422856
    "J10" = 0.1036893349422857
    "M10" = 11.01658666666667
    "N10" = 33.04976
    "O10" = 0.1837409300120545
    "P10" = 0.1837409300120545
    "Q10" = 159.13918831664
    "R10" = 1432.25269484976
    "S10" = 0.01905197483462698
    "T10" = 0.01905197483462699
    "G11" = 14.445417
    "H11" = 43.336251
    "I11" = 0.1036893349422856
    "J11" = 0.1036893349422857
    "O11" = 0.06509859443116503
    "P11" = 0.06509859443116503
    "Q11" = 56.38230674923699
    "R11" = 507.4407607431329
    "S11" = 0.006750029962245082
    "T11" = 0.006750029962245082
    "G12" = 14.445417
    "H12" = 43.336251
    "I12" = 0.1036893349422856
    "J12" = 0.1036893349422857
    "M12" = 4.001997666666667
    "N12" = 12.005993
    "O12" = 0.06674760480978428
    "P12" = 0.06674760480978428
    "Q12" = 57.810525128027
    "R12" = 520.294726152243
    "S12" = 0.006921014751717039
    "T12" = 0.00692101475171704
    "G13" = 14.445417
    "H13" = 43.336251
    "I13" = 0.1036893349422856
    "J13" = 0.1036893349422857
    "M13" = 41.03546066666667
    "N13" = 123.106382
    "O13" = 0.6844128707469962
    "P13" = 0.6844128707469963
    "Q13" = 592.7743411170979
    "R13" = 5334.969070053881
    "S13" = 0.07096631539369654
    "T13" = 0.07096631539369656
    "G14" = 14.34140633333333
    "H14" = 43.024219
    "I14" = 0.1029427454285617
    "J14" = 0.1029427454285617
    "M14" = 11.01658666666667
    "N14" = 33.04976
    "O14" = 0.1837409300120545
    "P14" = 0.1837409300120545
    "Q14" = 157.9933457930489
    "R14" = 1421.94011213744
    "S14" = 0.01891479578303809
    "T14" = 0.0189147957830381
    "G15" = 14.34140633333333
    "H15" = 43.024219
    "I15" = 0.1029427454285617
    "J15" = 0.1029427454285617
    "O15" = 0.06509859443116503
    "P15" = 0.06509859443116503
    "Q15" = 55.97633983854188
    "R15" = 503.787058546877
    "S15" = 0.006701428034284604
    "T15" = 0.006701428034284604
    "G16" = 14.34140633333333
    "H16" = 43.024219
    "I16" = 0.1029427454285617
    "J16" = 0.1029427454285617
    "M16" = 4.001997666666667
    "N16" = 12.005993
    "O16" = 0.06674760480978428
    "P16" = 0.06674760480978428
    "Q16" = 57.39427468271856
    "R16" = 516.548472144467
    "S16" = 0.006871181689899862
    "T16" = 0.006871181689899863
    "G17" = 14.34140633333333
    "H17" = 43.024219
    "I17" = 0.1029427454285617
    "J17" = 0.1029427454285617
    "M17" = 41.03546066666667
    "N17" = 123.106382
    "O17" = 0.6844128707469962
    "P17" = 0.6844128707469963
    "Q17" = 588.5062154961842
    "R17" = 5296.555939465658
    "S17" = 0.07045533992133911
    "T17" = 0.07045533992133912

}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
